# Generate Report for Handoff
# Swaps the handed-off markdown file's GUID-named artifacts (and their
# handoff timestamps) for a fresh handoff run:
#   51ee8fdb-77d0-491e-8cf5-360cee01e55d  ->  9e604063-8613-4655-8e1d-1237b36e4ccf
# across the Overview / zh-cn / de-de sheets, including the hyperlink
# "display" text that mirrors each cell's visible value.

$wb = $excel.ActiveWorkbook

$oldGuid = "51ee8fdb-77d0-491e-8cf5-360cee01e55d"
$newGuid = "9e604063-8613-4655-8e1d-1237b36e4ccf"

$oldZhHash = "250a896602b8adeec6da0b7124a6c96c58e0e744"
$newZhHash = "6ec7d91b8c6dcda04a62f63769ccb01233daed97"

# -------------------------------------------------------------------------
# Overview sheet: A2 (source .md display/hyperlink) + D2 (latest handoff dt)
# -------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newGuid.md"
$wsOverview.Range("D2").Value2 = "2016-03-20 05:02:50"

# Note: Range.Hyperlinks.Delete() clears *every* hyperlink on the sheet
# (it is not scoped to the range), so gather all hyperlinks first, delete
# once, then re-add each one (with the updated display text where needed).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d6bffed709a1835d80526ed026d1a22f0eb45211/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)

# -------------------------------------------------------------------------
# zh-cn sheet: A2 (source .md) + D2 (handoff .xlf) + E2 (handoff datetime)
# -------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = "$newGuid.md"
$wsZhCn.Range("D2").Value2 = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value2 = "2016-03-20 05:02:42"

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d6bffed709a1835d80526ed026d1a22f0eb45211/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f93907a1973acf05f2c8d0fa26c6443535af265d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$newGuid.$newZhHash.zh-cn.xlf",
    "",
    "",
    "$newGuid.$newZhHash.zh-cn.xlf"
)

# -------------------------------------------------------------------------
# de-de sheet: A2 (source .md) + D2 (handoff .xlf) + E2 (handback datetime)
# -------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = "$newGuid.md"
$wsDeDe.Range("D2").Value2 = "$newGuid.$newZhHash.de-de.xlf"
$wsDeDe.Range("E2").Value2 = "2016-03-20 05:02:50"

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d6bffed709a1835d80526ed026d1a22f0eb45211/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b6e22e490695d12ccca3d88bdf225e9208d98ff7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$newGuid.$newZhHash.de-de.xlf",
    "",
    "",
    "$newGuid.$newZhHash.de-de.xlf"
)
